$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Programs to include" (index 2): clear the "x" marker from most rows,
# keep it on two rows, and move the selection.
# ---------------------------------------------------------------------------
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Activate()

$rowsToClear = @(2) + (4..38) + (41..48)
foreach ($r in $rowsToClear) {
    $wsInclude.Range("B$r").Value = $null
}

$wsInclude.Range("B41").Select()

# ---------------------------------------------------------------------------
# Sheet "Program dependencies" (index 5): move the selection.
# ---------------------------------------------------------------------------
$wsDeps = $wb.Worksheets.Item(5)
$wsDeps.Activate()
$wsDeps.Range("A39").Select()

# ---------------------------------------------------------------------------
# Sheet "Budget scenario" (index 8): move the selection.
# ---------------------------------------------------------------------------
$wsBudget = $wb.Worksheets.Item(8)
$wsBudget.Activate()
$wsBudget.Range("A44").Select()

# ---------------------------------------------------------------------------
# Sheet "Coverage scenario" (index 7): add three new year columns, clear a
# value, add two new values, change zoom and selection.
# ---------------------------------------------------------------------------
$wsCoverage = $wb.Worksheets.Item(7)
$wsCoverage.Activate()

$wsCoverage.Range("K1").Copy()
$wsCoverage.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsCoverage.Range("L1").Value = 2026
$wsCoverage.Range("M1").Value = 2027
$wsCoverage.Range("N1").Value = 2028

$wsCoverage.Range("E2").Value = $null
$wsCoverage.Range("D42").Value = 0.5
$wsCoverage.Range("D43").Value = 0.5

$excel.ActiveWindow.Zoom = 94
$wsCoverage.Range("D44").Select()

# ---------------------------------------------------------------------------
# Sheet "Programs cost and coverage" (index 3): change a few values, move the
# selection, and make it the active (tabSelected) sheet, activated last.
# ---------------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item(3)
$wsCost.Activate()

$wsCost.Range("B40").Value = 0.5
$wsCost.Range("B43").Value = 0
$wsCost.Range("B44").Value = 0

$wsCost.Range("B43").Select()

# ---------------------------------------------------------------------------
# Window placement / first visible tab (best effort).
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = -26020
$win.Top = -20940
$win.Width = 19200
$win.Height = 21140
